# Bond dates update: TODAY() advanced by one day (2023-10-03 -> 2023-10-04).
# Columns G ("Dni od poprzedniej wypłaty" = days since previous payment) and
# I ("Dni do następnej wypłaty" = days to next payment) are stored as plain
# computed numbers (no formulas), so every cell with a value in column G is
# incremented by 1 and every cell with a value in column I is decremented by 1.
# Cells that are empty (bonds with no previous payment date) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $gCell = $ws.Cells.Item($r, 7)
    $gVal = $gCell.Value2
    if ($gVal -ne $null) {
        $gCell.Value = $gVal + 1
    }

    $iCell = $ws.Cells.Item($r, 9)
    $iVal = $iCell.Value2
    if ($iVal -ne $null) {
        $iCell.Value = $iVal - 1
    }
}
